$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the Excel table by one row so the table ref / autoFilter / dimension grow to F71
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null

# Row 70: fill in the start/end time data that was previously missing
$ws.Range("A70").Value = 43394
$ws.Range("B70").Value = 0.3840277777777778
$ws.Range("C70").Value = 0.69652777777777775

# Row 71: new daily record (date only; start/end time left blank)
$ws.Range("A71").Value = 43395
$ws.Range("D71").Formula = "=(C71-B71)* 1440"
$ws.Range("E71").Formula = "=IF(C71>B71, (C71-B71)*1440, (B71-C71)*1440)"
$ws.Range("F71").Formula = "=ABS((C71-B71)*1440)"

# Update the active selection to reflect the newly-added row
$ws.Range("B71").Select()
